$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 342.30768
$ws.Range("I15").Value = 342.30768
$ws.Range("K15").Value = 1026.92304
$ws.Range("M15").Value = -857.9230400000001

$ws.Range("H33").Value = 436.32
$ws.Range("I33").Value = 296.1579
$ws.Range("K33").Value = 296.1579
$ws.Range("M33").Value = -67.15789999999998

$ws.Range("H58").Value = 907.3333
$ws.Range("I58").Value = 444
$ws.Range("K58").Value = 1332
$ws.Range("M58").Value = -1182

$ws.Range("H70").Value = 85506.75
$ws.Range("I70").Value = 143430.64
$ws.Range("J70").Value = 4413.3
$ws.Range("K70").Value = 430291.92
$ws.Range("L70").Value = 13239.9
$ws.Range("M70").Value = -430021.92
$ws.Range("N70").Value = -13779.9

$ws.Range("H73").Value = 85506.75
$ws.Range("I73").Value = 143430.64
$ws.Range("J73").Value = 4413.3
$ws.Range("K73").Value = 430291.92
$ws.Range("L73").Value = 13239.9
$ws.Range("M73").Value = -429355.92
$ws.Range("N73").Value = -15111.9

$ws.Range("H101").Value = 12711.5
$ws.Range("I101").Value = 11613
$ws.Range("J101").Value = 14542.333
$ws.Range("K101").Value = 34839
$ws.Range("L101").Value = 43626.999
$ws.Range("M101").Value = -33217
$ws.Range("N101").Value = -46870.999

$ws.Range("H108").Value = 79533.336
$ws.Range("J108").Value = 79533.336
$ws.Range("L108").Value = 79533.336
$ws.Range("N108").Value = -87213.336

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 2580.3333
$ws.Range("I132").Value = 2633.0908
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7899.2724
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5369.2724
$ws.Range("N132").Value = -11060

$ws.Range("H137").Value = 2918.625
$ws.Range("J137").Value = 4279.5454
$ws.Range("L137").Value = 12838.6362
$ws.Range("N137").Value = -17938.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1576.5
$ws.Range("I2").Value = 1569.1428
$ws.Range("K2").Value = 1569.1428
$ws.Range("M2").Value = -1456.1428

$ws.Range("H32").Value = 4565
$ws.Range("I32").Value = 2589.4866
$ws.Range("K32").Value = 2589.4866
$ws.Range("M32").Value = -2302.4866

$ws.Range("H45").Value = 2721.5557
$ws.Range("I45").Value = 965.6667
$ws.Range("K45").Value = 965.6667
$ws.Range("M45").Value = -588.6667

$ws.Range("H102").Value = 3206.889
$ws.Range("I102").Value = 3206.889
$ws.Range("K102").Value = 3206.889
$ws.Range("M102").Value = -1584.889

$ws.Range("H116").Value = 1576.5
$ws.Range("I116").Value = 1569.1428
$ws.Range("K116").Value = 1569.1428
$ws.Range("M116").Value = 724.8571999999999

$ws.Range("H132").Value = 3144.75
$ws.Range("I132").Value = 2376.6667
$ws.Range("K132").Value = 7130.000100000001
$ws.Range("M132").Value = -4600.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1576.5
$ws.Range("I3").Value = 1569.1428
$ws.Range("K3").Value = 1569.1428
$ws.Range("M3").Value = -1455.1428

$ws.Range("H86").Value = 2854.4
$ws.Range("I86").Value = 2758
$ws.Range("J86").Value = 2999
$ws.Range("K86").Value = 2758
$ws.Range("L86").Value = 2999
$ws.Range("M86").Value = -1635
$ws.Range("N86").Value = -5245

$ws.Range("H89").Value = 2854.4
$ws.Range("I89").Value = 2758
$ws.Range("J89").Value = 2999
$ws.Range("K89").Value = 13790
$ws.Range("L89").Value = 14995
$ws.Range("M89").Value = -8174
$ws.Range("N89").Value = -26227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2335.4546
$ws.Range("I16").Value = 2277.6667
$ws.Range("J16").Value = 2595.5
$ws.Range("K16").Value = 2277.6667
$ws.Range("L16").Value = 2595.5
$ws.Range("M16").Value = -1990.6667
$ws.Range("N16").Value = -3169.5

$ws.Range("H104").Value = 44999.75
$ws.Range("J104").Value = 44999.75
$ws.Range("L104").Value = 44999.75
$ws.Range("N104").Value = -50241.75

$ws.Range("H109").Value = 36666.332
$ws.Range("J109").Value = 36666.332
$ws.Range("L109").Value = 36666.332
$ws.Range("N109").Value = -38746.332

$ws.Range("H113").Value = 2335.4546
$ws.Range("I113").Value = 2277.6667
$ws.Range("J113").Value = 2595.5
$ws.Range("K113").Value = 2277.6667
$ws.Range("L113").Value = 2595.5
$ws.Range("M113").Value = -107.6667000000002
$ws.Range("N113").Value = -6935.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3834.5293
$ws.Range("I12").Value = 3961.25
$ws.Range("J12").Value = 3721.889
$ws.Range("K12").Value = 11883.75
$ws.Range("L12").Value = 11165.667
$ws.Range("M12").Value = -11710.75
$ws.Range("N12").Value = -11511.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H80").Value = 15813.223
$ws.Range("I80").Value = 4463.2
$ws.Range("J80").Value = 30000.75
$ws.Range("K80").Value = 4463.2
$ws.Range("L80").Value = 30000.75
$ws.Range("M80").Value = -3465.2
$ws.Range("N80").Value = -31996.75

$ws.Range("H83").Value = 15813.223
$ws.Range("I83").Value = 4463.2
$ws.Range("J83").Value = 30000.75
$ws.Range("K83").Value = 22316
$ws.Range("L83").Value = 150003.75
$ws.Range("M83").Value = -17324
$ws.Range("N83").Value = -159987.75

$ws.Range("H107").Value = 607.1177
$ws.Range("I107").Value = 157.45454
$ws.Range("J107").Value = 1431.5
$ws.Range("K107").Value = 157.45454
$ws.Range("L107").Value = 1431.5
$ws.Range("M107").Value = 1762.54546
$ws.Range("N107").Value = -5271.5

$ws.Range("H122").Value = 2600.4167
$ws.Range("I122").Value = 2673.4546
$ws.Range("K122").Value = 8020.3638
$ws.Range("M122").Value = -5570.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3184.4666
$ws.Range("I46").Value = 2446.375
$ws.Range("J46").Value = 4028
$ws.Range("K46").Value = 2446.375
$ws.Range("L46").Value = 4028
$ws.Range("M46").Value = -2258.375
$ws.Range("N46").Value = -4404

$ws.Range("H61").Value = 3612.5
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H68").Value = 3566.5557
$ws.Range("I68").Value = 3012
$ws.Range("K68").Value = 3012
$ws.Range("M68").Value = -2263

$ws.Range("H71").Value = 3566.5557
$ws.Range("I71").Value = 3012
$ws.Range("K71").Value = 15060
$ws.Range("M71").Value = -11316

$ws.Range("H82").Value = 2999.5557
$ws.Range("J82").Value = 2999
$ws.Range("L82").Value = 2999
$ws.Range("N82").Value = -3721

$ws.Range("H85").Value = 2999.5557
$ws.Range("J85").Value = 2999
$ws.Range("L85").Value = 2999
$ws.Range("N85").Value = -5495

$ws.Range("H113").Value = 3612.5
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 4878.6875
$ws.Range("I122").Value = 4818.6665
$ws.Range("J122").Value = 4955.857
$ws.Range("K122").Value = 14455.9995
$ws.Range("L122").Value = 14867.571
$ws.Range("M122").Value = -12005.9995
$ws.Range("N122").Value = -19767.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H70").Value = 26547.5
$ws.Range("I70").Value = 30095
$ws.Range("J70").Value = 23000
$ws.Range("K70").Value = 30095
$ws.Range("L70").Value = 23000
$ws.Range("M70").Value = -29780
$ws.Range("N70").Value = -23630

$ws.Range("H73").Value = 26547.5
$ws.Range("I73").Value = 30095
$ws.Range("J73").Value = 23000
$ws.Range("K73").Value = 30095
$ws.Range("L73").Value = 23000
$ws.Range("M73").Value = -29003
$ws.Range("N73").Value = -25184

$ws.Range("H108").Value = 91999.5
$ws.Range("J108").Value = 91999.5
$ws.Range("L108").Value = 91999.5
$ws.Range("N108").Value = -99679.5

$ws.Range("H113").Value = 450.25
$ws.Range("I113").Value = 450.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1350.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 819.25
$ws.Range("N113").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H121").Value = 109999
$ws.Range("J121").Value = 109999
$ws.Range("L121").Value = 109999
$ws.Range("N121").Value = -113493

$ws.Range("H136").Value = 3830.2307
$ws.Range("I136").Value = 1979.8
$ws.Range("J136").Value = 9998.333000000001
$ws.Range("K136").Value = 5939.4
$ws.Range("L136").Value = 29994.999
$ws.Range("M136").Value = -3389.4
$ws.Range("N136").Value = -35094.999

$ws.Range("H141").Value = 98000
$ws.Range("J141").Value = 98000
$ws.Range("L141").Value = 98000
$ws.Range("N141").Value = -108360
